$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (existing rows 2-14 shift down to 3-15)
$ws.Rows.Item(2).Insert()
# Drop the formatting the new row inherited from the header row above
$ws.Rows.Item(2).ClearFormats()

$r = 2

function Set-TextCell($cell, $val) {
    # Force text storage so date-looking strings ("2024-05-13", ...) are not
    # auto-converted into date serial numbers by the smart-entry parser.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $ws.Cells.Item($r, 1) "2024-05-13"
$ws.Cells.Item($r, 2).Value = "노브랜드"
$ws.Cells.Item($r, 3).Value = "삼성"
Set-TextCell $ws.Cells.Item($r, 4) "2024-05-17"
Set-TextCell $ws.Cells.Item($r, 5) "2024-05-23"
$ws.Cells.Item($r, 6).Value = 16800000
$ws.Cells.Item($r, 7).Value = 1200000
$ws.Cells.Item($r, 8).Value = "-"
$ws.Cells.Item($r, 9).Value = 8700
$ws.Cells.Item($r, 10).Value = 11000
$ws.Cells.Item($r, 11).Value = "-"
$ws.Cells.Item($r, 12).Value = 14000
$ws.Cells.Item($r, 13).Value = "-"
$ws.Cells.Item($r, 14).Value = "-"
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = "-"
$ws.Cells.Item($r, 17).Value = "-"
$ws.Cells.Item($r, 18).Value = "2071 : 1"
$ws.Cells.Item($r, 19).Value = "-"
$ws.Cells.Item($r, 20).Value = "-"
